$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26

$ws.Range("A" + $row).Value = 'remas ali almadani_20251202_134130'
$ws.Range("C" + $row).Value = 'remas ali almadani'
$ws.Range("D" + $row).Value = 19
$ws.Range("E" + $row).Value = 'Female'
$ws.Range("F" + $row).Value = '2025-12-02 13:41:30'
$ws.Range("G" + $row).Value = '{
  "portion": 0.2,
  "diet": 0.2857142857142857,
  "salt": 0.6,
  "fat": 0.6,
  "natural": 0.8,
  "convenience": 0.4,
  "price": 1.0
}'
$ws.Range("H" + $row).Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range("I" + $row).Value = "'" + '0.578'
$ws.Range("J" + $row).Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("K" + $row).Value = 'Maruchan Ramen Sabor Pollo'
$ws.Range("L" + $row).Value = "'" + '0.566'
$ws.Range("M" + $row).Value = 'Sabor clásico, económico, alto en sodio, no saludable, nostálgico'
$ws.Range("N" + $row).Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range("O" + $row).Value = "'" + '0.455'
$ws.Range("P" + $row).Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("Q" + $row).Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range("R" + $row).Value = "'" + '0.712'
$ws.Range("S" + $row).Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("T" + $row).Value = 'Annie’s Shells & White Cheddar'
$ws.Range("U" + $row).Value = "'" + '0.625'
$ws.Range("V" + $row).Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range("W" + $row).Value = 'Velveeta Original Shells & Cheese (microwave cups)'
$ws.Range("X" + $row).Value = "'" + '0.567'
$ws.Range("Y" + $row).Value = 'Muy cremoso, porción individual, rápido, salado, ideal para niños'
$ws.Range("Z" + $row).Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range("AA" + $row).Value = "'" + '0.657'
$ws.Range("AB" + $row).Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range("AC" + $row).Value = 'Jack Link’s Beef Jerky Original'
$ws.Range("AD" + $row).Value = "'" + '0.656'
$ws.Range("AE" + $row).Value = 'Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña'
$ws.Range("AF" + $row).Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range("AG" + $row).Value = "'" + '0.644'
$ws.Range("AH" + $row).Value = 'Portátil, saludable, fácil, buena textura, sabor suave'

# Reset row height back to standard after the multi-line JSON text in G26
# triggers an auto row-height bump, so the serialized row matches the
# no-custom-height row produced by the source tooling.
$ws.Rows.Item($row).AutoFit()
